# Fill in the second row of each of the three temperature tables:
# update the "To" date and populate Average / Max / Min columns.

$d = $word.ActiveDocument

# Values per table, in document order: Wroclaw, Warsaw, Berlin
$values = @(
    @{ To = "06.12.2016"; Average = "12"; Max = "13"; Min = "10" },
    @{ To = "06.12.2016"; Average = "12"; Max = "15"; Min = "8"  },
    @{ To = "06.12.2016"; Average = "10"; Max = "14"; Min = "7"  }
)

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $table = $d.Tables.Item($i)
    $v = $values[$i - 1]

    # Row 2 holds the data: From, To, Average, Max, Min
    $table.Cell(2, 2).Range.Text = $v.To
    $table.Cell(2, 3).Range.Text = $v.Average
    $table.Cell(2, 4).Range.Text = $v.Max
    $table.Cell(2, 5).Range.Text = $v.Min
}
